# Applies the schedule update for professor Cláudio (ELM-2NA now in use).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 (20:00)
$ws.Range("B18").Value = "['ELM-2NA-Tecnologia da Soldagem', -, -, -]"
$ws.Range("C18").Value = "['MEC-2NB-Soldagem', -, 'ELM-2NA-Tecnologia da Soldagem', -]"
$ws.Range("F18").Value = "-"

# Row 19 (20:50)
$ws.Range("C19").Value = "['MEC-2NB-Soldagem', -, -, -]"
$ws.Range("F19").Value = "-"

# Row 20 (21:40)
$ws.Range("C20").Value = "['MEC-2NB-Soldagem', -, 'ELM-2NA-Tecnologia da Soldagem', -]"
$ws.Range("E20").Value = "[-, -, 'MEC-2NB-Soldagem', -]"
$ws.Range("F20").Value = "-"

# Row 21 (22:35)
$ws.Range("B21").Value = "[-, -, 'ELM-2NA-Tecnologia da Soldagem', -]"
$ws.Range("F21").Value = "-"
